$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 20
$ws.Range("H20").Value = 220
$ws.Range("I20").Value = 220
$ws.Range("K20").Value = 220
$ws.Range("M20").Value = 10
# Row 28
$ws.Range("H28").Value = 6982.4585
$ws.Range("I28").Value = 301.7857
$ws.Range("K28").Value = 301.7857
$ws.Range("M28").Value = 183.2143
# Row 32
$ws.Range("H32").Value = 12949.857
$ws.Range("J32").Value = 19500.25
$ws.Range("L32").Value = 19500.25
$ws.Range("N32").Value = -20152.25
# Row 35
$ws.Range("H35").Value = 220
$ws.Range("I35").Value = 220
$ws.Range("K35").Value = 220
$ws.Range("M35").Value = 159
# Row 51
$ws.Range("H51").Value = 19611.666
$ws.Range("I51").Value = 42326.332
$ws.Range("J51").Value = 13933
$ws.Range("K51").Value = 42326.332
$ws.Range("L51").Value = 13933
$ws.Range("M51").Value = -41842.332
$ws.Range("N51").Value = -14901
# Row 107
$ws.Range("H107").Value = 7486.077
$ws.Range("I107").Value = 7528.25
$ws.Range("K107").Value = 7528.25
$ws.Range("M107").Value = -5608.25
# Row 113
$ws.Range("H113").Value = 2610.2778
$ws.Range("I113").Value = 2599.125
$ws.Range("K113").Value = 2599.125
$ws.Range("M113").Value = 654.875
# Row 138
$ws.Range("H138").Value = 4412.661
$ws.Range("J138").Value = 4541.5884
$ws.Range("L138").Value = 13624.7652
$ws.Range("N138").Value = -23904.7652

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 5793.4165
$ws.Range("I2").Value = 2303.4285
$ws.Range("J2").Value = 10679.4
$ws.Range("K2").Value = 2303.4285
$ws.Range("L2").Value = 10679.4
$ws.Range("M2").Value = -2190.4285
$ws.Range("N2").Value = -10905.4
# Row 32
$ws.Range("H32").Value = 10092.75
$ws.Range("I32").Value = 6170.9873
$ws.Range("J32").Value = 23997.182
$ws.Range("K32").Value = 6170.9873
$ws.Range("L32").Value = 23997.182
$ws.Range("M32").Value = -5883.9873
$ws.Range("N32").Value = -24571.182
# Row 45
$ws.Range("H45").Value = 2623.2632
$ws.Range("I45").Value = 2292.3572
$ws.Range("K45").Value = 2292.3572
$ws.Range("M45").Value = -1915.3572
# Row 61
$ws.Range("H61").Value = 5000.3105
$ws.Range("I61").Value = 5021.75
$ws.Range("J61").Value = 4952.6665
$ws.Range("K61").Value = 5021.75
$ws.Range("L61").Value = 4952.6665
$ws.Range("M61").Value = -4809.75
$ws.Range("N61").Value = -5376.6665
# Row 63
$ws.Range("H63").Value = 3706.75
$ws.Range("J63").Value = 5349.8
$ws.Range("L63").Value = 5349.8
$ws.Range("N63").Value = -6721.8
# Row 66
$ws.Range("H66").Value = 3706.75
$ws.Range("J66").Value = 5349.8
$ws.Range("L66").Value = 26749
$ws.Range("N66").Value = -33613
# Row 116
$ws.Range("H116").Value = 5793.4165
$ws.Range("I116").Value = 2303.4285
$ws.Range("J116").Value = 10679.4
$ws.Range("K116").Value = 2303.4285
$ws.Range("L116").Value = 10679.4
$ws.Range("M116").Value = -9.428499999999985
$ws.Range("N116").Value = -15267.4
# Row 122
$ws.Range("H122").Value = 8211.4
$ws.Range("I122").Value = 7012.1763
$ws.Range("J122").Value = 10759.75
$ws.Range("K122").Value = 21036.5289
$ws.Range("L122").Value = 32279.25
$ws.Range("M122").Value = -18586.5289
$ws.Range("N122").Value = -37179.25
# Row 136
$ws.Range("H136").Value = 5000.3105
$ws.Range("I136").Value = 5021.75
$ws.Range("J136").Value = 4952.6665
$ws.Range("K136").Value = 15065.25
$ws.Range("L136").Value = 14857.9995
$ws.Range("M136").Value = -12515.25
$ws.Range("N136").Value = -19957.9995

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 5793.4165
$ws.Range("I3").Value = 2303.4285
$ws.Range("J3").Value = 10679.4
$ws.Range("K3").Value = 2303.4285
$ws.Range("L3").Value = 10679.4
$ws.Range("M3").Value = -2189.4285
$ws.Range("N3").Value = -10907.4
# Row 20
$ws.Range("H20").Value = 11907826
$ws.Range("I20").Value = 28573388
$ws.Range("J20").Value = 3852
$ws.Range("K20").Value = 28573388
$ws.Range("L20").Value = 3852
$ws.Range("M20").Value = -28573141
$ws.Range("N20").Value = -4346
# Row 54
$ws.Range("H54").Value = 4011.25
$ws.Range("I54").Value = 4011.25
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 4011.25
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -3527.25
$ws.Range("N54").ClearContents()
# Row 88
$ws.Range("H88").Value = 11565.571
$ws.Range("J88").Value = 12159.833
$ws.Range("L88").Value = 12159.833
$ws.Range("N88").Value = -12971.833
# Row 91
$ws.Range("H91").Value = 11565.571
$ws.Range("J91").Value = 12159.833
$ws.Range("L91").Value = 12159.833
$ws.Range("N91").Value = -14967.833
# Row 99
$ws.Range("H99").Value = 10640.218
$ws.Range("I99").Value = 11036.784
$ws.Range("J99").Value = 10181.6875
$ws.Range("K99").Value = 11036.784
$ws.Range("L99").Value = 10181.6875
$ws.Range("M99").Value = -9538.784
$ws.Range("N99").Value = -13177.6875
# Row 137
$ws.Range("H137").Value = 49000
$ws.Range("J137").Value = 49000
$ws.Range("L137").Value = 49000
$ws.Range("N137").Value = -59200

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 86
$ws.Range("H86").Value = 9274
$ws.Range("I86").Value = 6565.467
$ws.Range("K86").Value = 6565.467
$ws.Range("M86").Value = -5442.467
# Row 89
$ws.Range("H89").Value = 9274
$ws.Range("I89").Value = 6565.467
$ws.Range("K89").Value = 32827.335
$ws.Range("M89").Value = -27211.335
# Row 122
$ws.Range("H122").Value = 6279
$ws.Range("I122").Value = 5873.5713
$ws.Range("J122").Value = 7698
$ws.Range("K122").Value = 17620.7139
$ws.Range("L122").Value = 23094
$ws.Range("M122").Value = -15170.7139
$ws.Range("N122").Value = -27994
# Row 141
$ws.Range("H141").Value = 161826.16
$ws.Range("J141").Value = 163887.03
$ws.Range("L141").Value = 163887.03
$ws.Range("N141").Value = -174247.03

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1580.2413
$ws.Range("I5").Value = 864.6667
$ws.Range("J5").Value = 2346.9285
$ws.Range("K5").Value = 2594.0001
$ws.Range("L5").Value = 7040.7855
$ws.Range("M5").Value = -2482.0001
$ws.Range("N5").Value = -7264.7855
# Row 37
$ws.Range("H37").Value = 116509.87
$ws.Range("J37").Value = 116509.87
$ws.Range("L37").Value = 349529.61
$ws.Range("N37").Value = -349753.61
# Row 99
$ws.Range("H99").Value = 7877.4
$ws.Range("I99").Value = 5121.75
$ws.Range("K99").Value = 15365.25
$ws.Range("M99").Value = -13119.25
# Row 102
$ws.Range("H102").Value = 19266.666
$ws.Range("I102").Value = 18000
$ws.Range("K102").Value = 54000
$ws.Range("M102").Value = -51566
# Row 135
$ws.Range("H135").Value = 1580.2413
$ws.Range("I135").Value = 864.6667
$ws.Range("J135").Value = 2346.9285
$ws.Range("K135").Value = 7782.0003
$ws.Range("L135").Value = 21122.3565
$ws.Range("M135").Value = -5247.0003
$ws.Range("N135").Value = -26192.3565

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 58832850
$ws.Range("I80").Value = 125006000
$ws.Range("J80").Value = 12273.444
$ws.Range("K80").Value = 125006000
$ws.Range("L80").Value = 12273.444
$ws.Range("M80").Value = -125005002
$ws.Range("N80").Value = -14269.444
# Row 83
$ws.Range("H83").Value = 58832850
$ws.Range("I83").Value = 125006000
$ws.Range("J83").Value = 12273.444
$ws.Range("K83").Value = 625030000
$ws.Range("L83").Value = 61367.22
$ws.Range("M83").Value = -625025008
$ws.Range("N83").Value = -71351.22
# Row 97
$ws.Range("H97").Value = 1077.8667
$ws.Range("I97").Value = 1058.88
$ws.Range("K97").Value = 1058.88
$ws.Range("M97").Value = -562.8800000000001
# Row 102
$ws.Range("H102").Value = 4069.9055
$ws.Range("I102").Value = 2979.2856
$ws.Range("J102").Value = 6207.52
$ws.Range("K102").Value = 2979.2856
$ws.Range("L102").Value = 6207.52
$ws.Range("M102").Value = -1357.2856
$ws.Range("N102").Value = -9451.52
# Row 113
$ws.Range("H113").Value = 7451.8423
$ws.Range("I113").Value = 7525
$ws.Range("J113").Value = 7432.3335
$ws.Range("K113").Value = 7525
$ws.Range("L113").Value = 7432.3335
$ws.Range("M113").Value = -5355
$ws.Range("N113").Value = -11772.3335
# Row 137
$ws.Range("H137").Value = 100747.5
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 6990
$ws.Range("I22").Value = 970
$ws.Range("K22").Value = 970
$ws.Range("M22").Value = -675
# Row 27
$ws.Range("H27").Value = 6990
$ws.Range("I27").Value = 970
$ws.Range("K27").Value = 970
$ws.Range("M27").Value = -863
# Row 40
$ws.Range("H40").Value = 5436.276
$ws.Range("I40").Value = 4715.9414
$ws.Range("J40").Value = 6456.75
$ws.Range("K40").Value = 4715.9414
$ws.Range("L40").Value = 6456.75
$ws.Range("M40").Value = -4579.9414
$ws.Range("N40").Value = -6728.75
# Row 122
$ws.Range("H122").Value = 963070.6
$ws.Range("I122").Value = 1186930
$ws.Range("K122").Value = 3560790
$ws.Range("M122").Value = -3558340
# Row 133
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -65060

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2533.1177
$ws.Range("I81").Value = 2029
$ws.Range("J81").Value = 3253.2856
$ws.Range("K81").Value = 4058
$ws.Range("L81").Value = 6506.5712
$ws.Range("M81").Value = -2997
$ws.Range("N81").Value = -8628.5712
# Row 84
$ws.Range("H84").Value = 2533.1177
$ws.Range("I84").Value = 2029
$ws.Range("J84").Value = 3253.2856
$ws.Range("K84").Value = 20290
$ws.Range("L84").Value = 32532.856
$ws.Range("M84").Value = -14986
$ws.Range("N84").Value = -43140.856
# Row 122
$ws.Range("H122").Value = 9715.509
$ws.Range("I122").Value = 3774
$ws.Range("K122").Value = 11322
$ws.Range("M122").Value = -8872
